# Konnect Bill Payment Verification Checks added
# Adds columns I:P (bene operation / to-account / account title / bene bank /
# bene id queries + db_val) to the BeneDeletion.xlsx data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Group 1: columns I (bene_op_type_query) and J (bene_op_value) ----
$ws.Range("I1").Value2 = "bene_op_type_query"
$ws.Range("J1").Value2 = "bene_op_value"
$ws.Range("I2").Value2 = "SELECT BENE_OPERATION_TYPE FROM DC_TRANSACTION DT where DT.TRANSACTION_ID='{TRANSACTION_ID}'"
$ws.Range("J2").Value2 = "REMOVE"

# ---- Group 2: columns K..O (to_account / account_title / bene_bank / bene_id_tran / bene_id) ----
$ws.Range("K1").Value2 = "to_account_query"
$ws.Range("L1").Value2 = "account_title_query"
$ws.Range("M1").Value2 = "bene_bank_query"
$ws.Range("N1").Value2 = "bene_id_tran_query"
$ws.Range("O1").Value2 = "bene_id_query"
$ws.Range("K2").Value2 = "SELECT K.TO_ACCOUNT FROM DC_TRANSACTION K WHERE K.TRANSACTION_ID = '{TRANSACTION_ID}'"
$ws.Range("L2").Value2 = "SELECT K.FT_TO_ACCOUNT_TITLE FROM DC_TRANSACTION K WHERE K.TRANSACTION_ID = '{TRANSACTION_ID}'"
$ws.Range("M2").Value2 = "SELECT K.BENEFICIARY_BANK FROM DC_TRANSACTION K WHERE K.TRANSACTION_ID = '{TRANSACTION_ID}'"
$ws.Range("N2").Value2 = "SELECT BENEFICIARY_ID FROM DC_TRANSACTION K WHERE K.TRANSACTION_ID = '{TRANSACTION_ID}'"
$ws.Range("O2").Value2 = "SELECT FUND_TRANSFER_BENEFICIARY_ID FROM DC_FUND_TRANSFER_BENEFICIARY K WHERE K.CUSTOMER_INFO_ID = (Select CUSTOMER_INFO_ID from DC_CUSTOMER_INFO L WHERE L.CUSTOMER_NAME = '{customer_name}' ) and K.ACCOUNT_NO = '{account_number}'"

# ---- Group 3: column P (db_val) ----
$ws.Range("P1").Value2 = "db_val"
$ws.Range("P2").Value2 = "DIGITAL_CHANNEL_SEC"

# ---- Row 3 repeats row 2's values (same shared strings, no new ones) ----
$ws.Range("I3").Value2 = $ws.Range("I2").Value2
$ws.Range("J3").Value2 = $ws.Range("J2").Value2
$ws.Range("K3").Value2 = $ws.Range("K2").Value2
$ws.Range("L3").Value2 = $ws.Range("L2").Value2
$ws.Range("M3").Value2 = $ws.Range("M2").Value2
$ws.Range("N3").Value2 = $ws.Range("N2").Value2
$ws.Range("O3").Value2 = $ws.Range("O2").Value2
$ws.Range("P3").Value2 = $ws.Range("P2").Value2

# ---- Number format: columns I,K use Text format "@" (style index 1,
# the same style already used by header/query columns A,B,C,E,F,G,H) ----
foreach ($addr in "I1", "J1", "K1", "L1", "M1", "P1", "I2", "K2", "I3", "K3") {
    $ws.Range($addr).NumberFormat = "@"
}

# ---- Font: columns N,O use a distinct font (Calibri, no explicit colour) ----
foreach ($addr in "N1", "O1", "N2", "O2", "N3", "O3") {
    $ws.Range($addr).Font.Name = "Calibri"
}

# ---- Column widths (best-fit, mirrors Excel's "AutoFit Column Width") ----
$ws.Columns.Item(9).AutoFit() | Out-Null
$ws.Columns.Item(10).AutoFit() | Out-Null
$ws.Columns.Item(11).AutoFit() | Out-Null
$ws.Columns.Item(12).AutoFit() | Out-Null
$ws.Columns.Item(13).AutoFit() | Out-Null
$ws.Columns.Item(14).AutoFit() | Out-Null
$ws.Columns.Item(15).AutoFit() | Out-Null
$ws.Columns.Item(16).AutoFit() | Out-Null

# ---- Selection moves to D18, as left by the author after editing ----
$ws.Range("D18").Select() | Out-Null
